$wb = $excel.ActiveWorkbook
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

$wsMeta.Range("B3").Value = "2.0.0"
$wsMeta.Range("B8").Value = "2025-10-20T13:10:23+00:00"
$wsMeta.Range("B12").Value = @"

- **Séjour** : commentaire relatif au séjour.
- **Événement** : commentaires sur le déroulé de l’évènement.
- **Évaluation** : commentaire libre sur le contenu ou le résultat de l’évaluation.
- **Champ évalué** : commentaire spécifique à un item ou sous-item évalué.
"@

$wsElem.Range("L6").Value = "Exemple de commentaire : Cet évènement a débuté plus tard l’usager était sous la douche à l’heure du début du rendez-vous."
